# SSEL-LR.xlsx edit: rename sheet and update the saved selection/scroll
# position, per the commit "Edicion Lista de Requitos / Mejorando lista
# de requisitos".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet tab from "Hoja 1" to "Lista de Requisitos".
$ws.Name = "Lista de Requisitos"

# 2. Make sure this sheet is the active/selected one, then move the
#    visible scroll window and the active cell/selection to match where
#    the author left off working (column E.. row 4 in view, H12 selected).
$ws.Activate()

$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 5
$win.Zoom = 70

$ws.Range("H12").Select() | Out-Null
